$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # "基金受益憑證" (fund) sheet

# --- Row 1 used to be a duplicate of row 2's data; turn it into a real
#     header row, and extend it with the new metadata columns (I1:O1) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# --- Rows 2-6: fill in the new metadata columns I:O for every existing
#     data row (property_category / category / date / legislator_name /
#     legislator_id / source_file / index) ---
$indexes = @(71, 72, 73, 74, 77)
for ($i = 0; $i -lt $indexes.Length; $i++) {
    $r = $i + 2

    $ws.Range("I$r").Value = "fund"
    $ws.Range("J$r").Value = "normal"

    # Force text so Excel doesn't auto-convert the ISO-ish date string into
    # a date serial number.
    $ws.Range("K$r").NumberFormat = "@"
    $ws.Range("K$r").Value = "2012-04-24"

    $ws.Range("L$r").Value = "徐耀昌"
    $ws.Range("M$r").Value = 921
    $ws.Range("N$r").Value = "tmp6e501"
    $ws.Range("O$r").Value = $indexes[$i]
}
